$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = -1
$ws.Range("C2").Value = 2.6666795265009045
$ws.Range("C3").Value = 4.3999737722116592

$ws.Columns.Item(3).ColumnWidth = 10.8333333333333
